# This "edit" originates from a commit whose actual payload is a refactor
# of the (unrelated) DSL/table-accessor tooling used to *generate* this
# test-fixture .docx. Re-running the generator with a newer docx4j/JVM on
# a different OS only changes incidental, non-semantic artifacts of that
# tool: the order in which XML namespace declarations are written on the
# root elements of word/document.xml, word/footer.xml and word/styles.xml,
# and the free-text "Created by docx4j ..." generator banner left as an
# XML comment inside word/document.xml's body.
#
# None of that is reachable (or meaningful) through Word's document object
# model: namespace-prefix ordering is decided by the XML serializer, not
# by any document property, and the generator banner is not Word body
# text, not a Word "Comment" (annotation), and not a custom XML part -
# it is inert tooling metadata that a real editing action in Word can
# neither see nor change. Actually touching the document body through the
# Word OM (even a no-op formatting change) would force Word to rewrite
# word/document.xml from its internal model and silently drop that
# foreign comment, which would be a much larger, unfaithful change than
# the one described by the diff.
#
# The diff therefore carries no reachable/semantic document edit for this
# file's actual content (same text, same runs, same formatting, same
# structure throughout). Intentionally perform no content mutation here.
$d = $word.ActiveDocument
